$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Title($ws, $addr) {
    $ws.Range($addr).Font.Bold = $true
}
function Set-Source($ws, $addr) {
    $ws.Range($addr).Font.Italic = $true
}

# Remove the existing hyperlink (it will be re-added at its new location below)
$ws.Range("A49").Hyperlinks.Delete()

# Clear the old "Sector Distribution Details" block so nothing stale is left behind
$ws.Range("A23:D54").Clear()

$ws.Range("B23").Value = "Number of employees"
Set-Title $ws "B23"
$ws.Range("C23").Value = "Assets (local currency, unless noted otherwise)"
Set-Title $ws "C23"
$ws.Range("D23").Value = "Turnover (local currency, unless noted otherwise)"
Set-Title $ws "D23"
$ws.Range("A24").Value = "Micro"
$ws.Range("B24").Value = "1-9"
$ws.Range("C24").Value = ""
$ws.Range("D24").Value = ""
$ws.Range("A25").Value = "Small"
$ws.Range("B25").Value = "10-49"
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = ""
$ws.Range("A26").Value = "Medium"
$ws.Range("B26").Value = "50-200"
$ws.Range("C26").Value = ""
$ws.Range("D26").Value = ""
$ws.Range("A27").Value = "Large"
$ws.Range("B27").Value = ">200"
$ws.Range("C27").Value = ""
$ws.Range("D27").Value = ""
$ws.Range("A30").Value = "Sector Distribution Details"
Set-Title $ws "A30"
$ws.Range("B32").Value = "MSMEs"
Set-Title $ws "B32"
$ws.Range("C32").Value = "%MSMEs"
Set-Title $ws "C32"
$ws.Range("A33").Value = "Manufacturing and Processing "
Set-Title $ws "A33"
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "19,413"
Set-Title $ws "B33"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "67.24"
Set-Title $ws "C33"
$ws.Range("A34").Value = "Cars, Motorcylces and goods for personal use"
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "16,512"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "57.19"
$ws.Range("A35").Value = "Other manufacturing"
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "2,901"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "10.05"
$ws.Range("A36").Value = "Services "
Set-Title $ws "A36"
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "6,358"
Set-Title $ws "B36"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "22.02"
Set-Title $ws "C36"
$ws.Range("A37").Value = "Construction"
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "276"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "0.96"
$ws.Range("A38").Value = "Health and Education"
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "324"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "1.12"
$ws.Range("A39").Value = "Hotels and Restaurants"
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "5,758"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "19.94"
$ws.Range("A40").Value = "Extractive Industries "
Set-Title $ws "A40"
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "53"
Set-Title $ws "B40"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "0.18"
Set-Title $ws "C40"
$ws.Range("A41").Value = "Infrastructure "
Set-Title $ws "A41"
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "377"
Set-Title $ws "B41"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "1.31"
Set-Title $ws "C41"
$ws.Range("A42").Value = "Electricity, gas and water"
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "31"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "0.11"
$ws.Range("A43").Value = "Transportation, storage and communications"
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "346"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "1.20"
$ws.Range("A44").Value = "Financial services "
Set-Title $ws "A44"
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "687"
Set-Title $ws "B44"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "2.38"
Set-Title $ws "C44"
$ws.Range("A45").Value = "Financial activities"
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "80"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "0.28"
$ws.Range("A46").Value = "Financial undertakings"
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "607"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "2.10"
$ws.Range("A47").Value = "Agribusiness "
Set-Title $ws "A47"
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "138"
Set-Title $ws "B47"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "0.48"
Set-Title $ws "C47"
$ws.Range("A48").Value = "Fishing"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "138"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "0.48"
$ws.Range("A49").Value = "Primary Agriculture "
Set-Title $ws "A49"
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "637"
Set-Title $ws "B49"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "2.21"
Set-Title $ws "C49"
$ws.Range("A50").Value = "Agriculture, animal production, hunting and forestry"
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "637"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "2.21"
$ws.Range("A51").Value = "Other "
Set-Title $ws "A51"
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "1,207"
Set-Title $ws "B51"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "4.18"
Set-Title $ws "C51"
$ws.Range("A52").Value = "Total "
Set-Title $ws "A52"
$ws.Range("B52").NumberFormat = "@"
$ws.Range("B52").Value = "28,870"
Set-Title $ws "B52"
$ws.Range("C52").NumberFormat = "@"
$ws.Range("C52").Value = "100"
Set-Title $ws "C52"
$ws.Range("A53").Value = "Source:"
Set-Source $ws "A53"
$ws.Range("A54").Value = "Instituto Nacional de Estatística de Moçambique (INE)"
Set-Source $ws "A54"
$ws.Range("A55").Value = "http://www.ine.gov.mz/censos_dir/cempre/resultadoscempre.pdf"
$ws.Range("A56").Value = "Page 6"
Set-Source $ws "A56"
$ws.Range("A59").Value = "AFDB"
Set-Title $ws "A59"
$ws.Range("A60").Value = "African Development Bank (AFDB), `"REPUBLIC OF MOZAMBIQUE: COUNTRY STRATEGY PAPER 2011-2015`", p. 5, 2011. Available at http://www.afdb.org/fileadmin/uploads/afdb/Documents/Policy-Documents/Mozambique%20-%202011-15%20CSP.pdf"
Set-Source $ws "A60"

# Re-create the hyperlink at its new location, reusing the same target URL
$ws.Hyperlinks.Add($ws.Range("A55"), "http://www.ine.gov.mz/censos_dir/cempre/resultadoscempre.pdf")

Write-Output "MSME Mozambique summary updated"
